$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.593.07"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "2.405.54"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.62"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.85"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.993"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.560"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").Value = "2.441.98"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0977"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.150"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("E13").Value = "  -4.88%  "
$ws.Range("D14").Value = "2.837.38"
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").Value = "57.442.59"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.93"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000134"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").Value = "2.427.26"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.30"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.13"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "314.55"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("E22").Value = "  +4.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("E24").Value = "  -1.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.67"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.992"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("D27").Value = "2.509.75"
$ws.Range("E27").Value = "  -2.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.385"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.153"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.52"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.16"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("D32").Value = "0.0₃0737"
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.21"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.15"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.992"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.17"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.23"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.86"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("E41").Value = "  +4.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.22"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.46"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "133.08"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +8.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.42"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.04"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "260.91"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.570"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0916"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0496"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0213"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.71%  "
